$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add new defined name (duplicate autofilter database range) scoped to the sheet
$ws.Names.Add("_xlnm._FilterDatabase_0_0", "=Requirement!`$A`$2:`$I`$12")

# 2. Populate the four new requirement rows (112-115)
# Row 112
$ws.Range("A112").Value = "UNIDEB_86"
$ws.Range("B112").Value = "R"
$ws.Range("C112").Value = "When the driver disables the TSR system it should reset it’s accumulated data."
$ws.Range("D112").Value = "New"
$ws.Range("E112").Value = "https://trello.com/c/sjrJZ7jM"
$ws.Range("F112").Value = "TSR"
$ws.Range("G112").Value = "Team1"
$ws.Range("H112").Value = "Sprint3"

# Row 113
$ws.Range("A113").Value = "UNIDEB_87"
$ws.Range("B113").Value = "R"
$ws.Range("C113").Value = "When a sign is enabled but the TSR system receives a sign that should cancel out the enabled sign the TSR system should send a disable signal via DON’T SHOW SUPPLEMENTAL SIGNS ON IC Signal (Signal id: 28)"
$ws.Range("D113").Value = "New"
$ws.Range("E113").Value = "https://trello.com/c/sjrJZ7jM"
$ws.Range("F113").Value = "TSR"
$ws.Range("G113").Value = "Team1"
$ws.Range("H113").Value = "Sprint3"

# Row 114
$ws.Range("A114").Value = "UNIDEB_88"
$ws.Range("B114").Value = "R"
$ws.Range("C114").Value = "The TSR should be able to disable certain signs based on elapsed time. In our case the “no speed limit” sign should get disables after a while."
$ws.Range("D114").Value = "New"
$ws.Range("E114").Value = "https://trello.com/c/sjrJZ7jM"
$ws.Range("F114").Value = "TSR"
$ws.Range("G114").Value = "Team1"
$ws.Range("H114").Value = "Sprint3"

# Row 115
$ws.Range("A115").Value = "UNIDEB_89"
$ws.Range("B115").Value = "R"
$ws.Range("C115").Value = "The cancel out signal for the speed limit sign should be the MOST RELEVANT SPEED LIMIT with data of 0."
$ws.Range("D115").Value = "New"
$ws.Range("E115").Value = "https://trello.com/c/sjrJZ7jM"
$ws.Range("F115").Value = "TSR"
$ws.Range("G115").Value = "Team1"
$ws.Range("H115").Value = "Sprint3"

# 3. Add hyperlinks for the trello links in column E, then restore their cell formatting
# (Hyperlinks.Add applies the builtin "Hyperlink" style; copy formatting back from column A of the same row)
$ws.Hyperlinks.Add($ws.Range("E112"), "https://trello.com/c/sjrJZ7jM", "", "", "https://trello.com/c/sjrJZ7jM")
$ws.Range("A112").Copy()
$ws.Range("E112").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("E113"), "https://trello.com/c/sjrJZ7jM", "", "", "https://trello.com/c/sjrJZ7jM")
$ws.Range("A113").Copy()
$ws.Range("E113").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("E114"), "https://trello.com/c/sjrJZ7jM", "", "", "https://trello.com/c/sjrJZ7jM")
$ws.Range("A114").Copy()
$ws.Range("E114").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("E115"), "https://trello.com/c/sjrJZ7jM", "", "", "https://trello.com/c/sjrJZ7jM")
$ws.Range("A115").Copy()
$ws.Range("E115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Update the active selection to mirror the author's final view
$ws.Range("C112").Select()
